$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 71-72, shifting existing rows 71-161 down to 73-163.
$ws.Rows("71:72").Insert()

# --- New row 71: Automn Giant ---
$ws.Cells.Item(71, 1).Value = 10
$ws.Cells.Item(71, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(71, 3).Value = "La Araucanía"
$ws.Cells.Item(71, 4).Value = 44579
$ws.Cells.Item(71, 5).Value = 9
$ws.Cells.Item(71, 6).Value = "Fruta"
$ws.Cells.Item(71, 7).Value = 100103
$ws.Cells.Item(71, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(71, 9).Value = 100103002
$ws.Cells.Item(71, 10).Value = "Ciruela"
$ws.Cells.Item(71, 11).Value = "Automn Giant"
$ws.Cells.Item(71, 12).Value = "Primera"
$ws.Cells.Item(71, 13).Value = 8
$ws.Cells.Item(71, 14).Value = 250000
$ws.Cells.Item(71, 15).Value = 250000
$ws.Cells.Item(71, 16).Value = 250000
$ws.Cells.Item(71, 17).Value = "`$/bins (420 kilos)"
$ws.Cells.Item(71, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(71, 19).Value = 595
$ws.Cells.Item(71, 20).Value = 420

# --- New row 72: Black Amber ---
$ws.Cells.Item(72, 1).Value = 10
$ws.Cells.Item(72, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(72, 3).Value = "La Araucanía"
$ws.Cells.Item(72, 4).Value = 44579
$ws.Cells.Item(72, 5).Value = 9
$ws.Cells.Item(72, 6).Value = "Fruta"
$ws.Cells.Item(72, 7).Value = 100103
$ws.Cells.Item(72, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(72, 9).Value = 100103002
$ws.Cells.Item(72, 10).Value = "Ciruela"
$ws.Cells.Item(72, 11).Value = "Black Amber"
$ws.Cells.Item(72, 12).Value = "Primera"
$ws.Cells.Item(72, 13).Value = 65
$ws.Cells.Item(72, 14).Value = 12000
$ws.Cells.Item(72, 15).Value = 12000
$ws.Cells.Item(72, 16).Value = 12000
$ws.Cells.Item(72, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(72, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(72, 19).Value = 667
$ws.Cells.Item(72, 20).Value = 18
